$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change column E (Obrigatorio) from "N" to "S" for rows 2 through 11
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("E$r").Value = "S"
}
